$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell address -> new value (SO number updated from 223 to 247 suffix)
$updates = [ordered]@{
    "A3"  = "JSSO1000247"
    "B3"  = "JSSO1000247"
    "C3"  = "JSSO1000247"
    "AJ3" = "JSCN1000247"
    "AL3" = "SLJSSO1000247"
    "AN3" = "MBLJSSO1000247"
    "AO3" = "HBLJSSO1000247"
}

# Scratch cell used to stash/restore original cell formatting (quote-prefix,
# fill, font, etc.) around the value change, since writing a new Value can
# otherwise reset formatting such as the text quote-prefix.
$helper = $ws.Range("ZZ1")

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)

    # Stash current formatting
    $rng.Copy()
    $helper.PasteSpecial(-4122)  # xlPasteFormats

    # Update the value
    $rng.Value = $updates[$addr]

    # Restore the original formatting
    $helper.Copy()
    $rng.PasteSpecial(-4122)  # xlPasteFormats
}

$helper.Clear()
$excel.CutCopyMode = $false
